$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Reformat the "dd/m/yyyy" style text dates (stored as literal text, not
#    real Excel dates) to "dd-mm-yyyy" with zero-padded day/month.
# ---------------------------------------------------------------------------
$dateUpdates = [ordered]@{
    "I4"  = "29-05-2020"
    "I5"  = "26-05-2020"
    "I6"  = "25-05-2020"
    "I7"  = "21-05-2020"
    "I8"  = "13-03-2020"
    "I9"  = "16-03-2020"
    "I10" = "18-02-2020"
    "I11" = "19-03-2020"
    "I12" = "17-03-2020"
    "I13" = "27-02-2020"
    "I14" = "20-03-2020"
    "I16" = "21-04-2020"
    "I17" = "27-04-2020"
    "I18" = "30-04-2020"
    "I19" = "30-04-2020"
    "I20" = "30-04-2020"
    "H22" = "13-06-2020"
    "H23" = "13-06-2020"
}

foreach ($ref in $dateUpdates.Keys) {
    $ws.Range($ref).Value = $dateUpdates[$ref]
}

# ---------------------------------------------------------------------------
# 2. H2 holds a real date serial (43957). Give it a new custom "long date"
#    number format, matching the newly-added numFmt / cellXf in styles.xml.
# ---------------------------------------------------------------------------
$ws.Range("H2").NumberFormat = '[$-F800]dddd\,\ mmmm\ dd\,\ yyyy'

# ---------------------------------------------------------------------------
# 3. Row heights: the rows holding the reformatted dates re-wrap slightly
#    differently once their text changed, so Excel recomputed their
#    auto-fit heights on save.
# ---------------------------------------------------------------------------
$rowHeights = [ordered]@{
    7  = 140.25
    8  = 120
    9  = 120
    11 = 120
    12 = 120
    13 = 120
    15 = 120
    16 = 120
    17 = 120
    18 = 120
    19 = 120
    20 = 120
    21 = 102
    22 = 120
    23 = 120
}

foreach ($r in $rowHeights.Keys) {
    $ws.Rows.Item($r).RowHeight = $rowHeights[$r]
}

# ---------------------------------------------------------------------------
# 4. View state: selection moved from H3 to H2, and the sheet scrolled so
#    row 21 / column F is the new top-left visible cell.
# ---------------------------------------------------------------------------
$ws.Range("H2").Select()
$excel.ActiveWindow.ScrollColumn = 6
$excel.ActiveWindow.ScrollRow = 21
